$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '26.280.87'

$ws.Cells.Item(3, 4).Value = '1.666.73'
$ws.Cells.Item(3, 5).Value = '  +0.72%  '

$ws.Cells.Item(4, 4).Value = "'1.010"
$ws.Cells.Item(4, 4).Style = 'Normal'
$ws.Cells.Item(4, 5).Value = '  +0.79%  '

$ws.Cells.Item(5, 4).Value = "'218.53"
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  +0.27%  '

$ws.Cells.Item(6, 4).Value = "'0.5326"
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  +1.47%  '

$ws.Cells.Item(7, 5).Value = '  +0.74%  '

$ws.Cells.Item(8, 4).Value = "'0.2640"
$ws.Cells.Item(8, 4).Style = 'Normal'
$ws.Cells.Item(8, 5).Value = '  +1.25%  '

$ws.Cells.Item(9, 4).Value = "'0.06381"
$ws.Cells.Item(9, 4).Style = 'Normal'
$ws.Cells.Item(9, 5).Value = '  +0.45%  '

$ws.Cells.Item(10, 4).Value = "'20.54"
$ws.Cells.Item(10, 4).Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  +0.78%  '

$ws.Cells.Item(11, 4).Value = "'0.07821"
$ws.Cells.Item(11, 4).Style = 'Normal'
$ws.Cells.Item(11, 5).Value = '  +0.22%  '

$ws.Cells.Item(12, 4).Value = "'4.566"
$ws.Cells.Item(12, 4).Style = 'Normal'
$ws.Cells.Item(12, 5).Value = '  +1.37%  '

$ws.Cells.Item(13, 4).Value = '1.668.39'
$ws.Cells.Item(13, 5).Value = '  -0.06%  '

$ws.Cells.Item(14, 4).Value = '1.894.27'
$ws.Cells.Item(14, 5).Value = '  +0.66%  '

$ws.Cells.Item(15, 4).Value = "'0.5534"
$ws.Cells.Item(15, 4).Style = 'Normal'
$ws.Cells.Item(15, 5).Value = '  +1.03%  '

$ws.Cells.Item(16, 4).Value = '0.0₅8203'
$ws.Cells.Item(16, 5).Value = '  -0.04%  '

$ws.Cells.Item(17, 4).Value = "'65.71"
$ws.Cells.Item(17, 4).Style = 'Normal'
$ws.Cells.Item(17, 5).Value = '  +0.47%  '

$ws.Cells.Item(18, 5).Value = '  +0.79%  '

$ws.Cells.Item(19, 4).Value = "'4.686"
$ws.Cells.Item(19, 4).Style = 'Normal'
$ws.Cells.Item(19, 5).Value = '  +2.40%  '

$ws.Cells.Item(20, 4).Value = "'194.02"
$ws.Cells.Item(20, 4).Style = 'Normal'
$ws.Cells.Item(20, 5).Value = '  +1.24%  '

$ws.Cells.Item(21, 4).Value = "'10.21"
$ws.Cells.Item(21, 4).Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  +1.51%  '

$ws.Cells.Item(22, 4).Value = "'6.034"
$ws.Cells.Item(22, 4).Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  +0.06%  '

$ws.Cells.Item(23, 4).Value = "'1.011"
$ws.Cells.Item(23, 4).Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  +0.75%  '

$ws.Cells.Item(24, 4).Value = "'145.61"
$ws.Cells.Item(24, 4).Style = 'Normal'
$ws.Cells.Item(24, 5).Value = '  +2.59%  '

$ws.Cells.Item(25, 4).Value = "'0.1228"
$ws.Cells.Item(25, 4).Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  -1.26%  '

$ws.Cells.Item(26, 4).Value = "'7.194"
$ws.Cells.Item(26, 4).Style = 'Normal'
$ws.Cells.Item(26, 5).Value = '  -0.89%  '

$ws.Cells.Item(27, 4).Value = "'16.12"
$ws.Cells.Item(27, 4).Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  -0.23%  '

$ws.Cells.Item(28, 4).Value = "'1.482"
$ws.Cells.Item(28, 4).Style = 'Normal'
$ws.Cells.Item(28, 5).Value = '  +3.67%  '

$ws.Cells.Item(29, 4).Value = "'0.05879"
$ws.Cells.Item(29, 4).Style = 'Normal'
$ws.Cells.Item(29, 5).Value = '  -0.62%  '

$ws.Cells.Item(30, 5).Value = '  +0.26%  '

$ws.Cells.Item(31, 4).Value = "'3.603"
$ws.Cells.Item(31, 4).Style = 'Normal'
$ws.Cells.Item(31, 5).Value = '  +2.22%  '

$ws.Cells.Item(32, 5).Value = '  +0.79%  '

$ws.Cells.Item(33, 5).Value = '  +1.40%  '

$ws.Cells.Item(34, 4).Value = "'0.9610"
$ws.Cells.Item(34, 4).Style = 'Normal'
$ws.Cells.Item(34, 5).Value = '  +0.81%  '

$ws.Cells.Item(35, 4).Value = "'2.825"
$ws.Cells.Item(35, 4).Style = 'Normal'
$ws.Cells.Item(35, 5).Value = '  +1.27%  '

$ws.Cells.Item(36, 5).Value = '  +0.49%  '

$ws.Cells.Item(37, 4).Value = "'0.5796"
$ws.Cells.Item(37, 4).Style = 'Normal'
$ws.Cells.Item(37, 5).Value = '  +1.71%  '

$ws.Cells.Item(38, 4).Value = "'0.01606"
$ws.Cells.Item(38, 4).Style = 'Normal'
$ws.Cells.Item(38, 5).Value = '  -0.78%  '

$ws.Cells.Item(39, 4).Value = "'0.8612"
$ws.Cells.Item(39, 4).Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  +1.35%  '

$ws.Cells.Item(40, 4).Value = "'5.845"
$ws.Cells.Item(40, 4).Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  +1.02%  '

$ws.Cells.Item(41, 2).Value = 'PaxDollar'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Cells.Item(41, 4).Value = "'1.009"
$ws.Cells.Item(41, 4).Style = 'Normal'
$ws.Cells.Item(41, 5).Value = '  +0.70%  '

$ws.Cells.Item(42, 2).Value = 'Maker'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(42, 4).Value = '1.048.70'
$ws.Cells.Item(42, 5).Value = '  +1.70%  '

$ws.Cells.Item(43, 4).Value = "'104.01"
$ws.Cells.Item(43, 4).Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  +0.86%  '

$ws.Cells.Item(44, 4).Value = '1.804.73'
$ws.Cells.Item(44, 5).Value = '  +0.38%  '

$ws.Cells.Item(45, 4).Value = "'57.68"
$ws.Cells.Item(45, 4).Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  +0.83%  '

$ws.Cells.Item(46, 4).Value = "'1.011"
$ws.Cells.Item(46, 4).Style = 'Normal'
$ws.Cells.Item(46, 5).Value = '  +1.27%  '

$ws.Cells.Item(47, 5).Value = '  -6.00%  '

$ws.Cells.Item(48, 4).Value = "'0.4379"
$ws.Cells.Item(48, 4).Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  +1.84%  '

$ws.Cells.Item(49, 4).Value = "'8.039"
$ws.Cells.Item(49, 4).Style = 'Normal'

$ws.Cells.Item(50, 4).Value = "'0.05158"
$ws.Cells.Item(50, 4).Style = 'Normal'
$ws.Cells.Item(50, 5).Value = '  -0.14%  '

$ws.Cells.Item(51, 5).Value = '  -3.10%  '
